$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the original (before) values for rows 2-8, columns D, M, N, O, P, Q, S
$cols = @("D","M","N","O","P","Q","S")
$snapshot = @{}
for ($r = 2; $r -le 8; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Mapping: target row -> source row (values come from the ORIGINAL/before data)
$mapping = @{
    2 = 3
    3 = 6
    4 = 5
    5 = 2
    6 = 7
    7 = 8
    8 = 4
}

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $src = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value2 = $src[$c]
    }
}
